$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.06692377132033253"
$ws.Range("E2").Value = [double]"0.06692377132033253"

$ws.Range("D4").Value = [double]"0.9999999728083765"
$ws.Range("E4").Value = [double]"0.9999999728083765"

$ws.Range("D5").Value = [double]"5.474639198905171E-102"
$ws.Range("E5").Value = [double]"5.474639198905171E-102"

$ws.Range("D6").Value = [double]"0.9999997924028738"
$ws.Range("E6").Value = [double]"0.9999997924028738"

$ws.Range("D8").Value = [double]"0.9999999557270978"
$ws.Range("E8").Value = [double]"4.427290223318892E-08"

$ws.Range("D9").Value = [double]"0.9999971239217129"
$ws.Range("E9").Value = [double]"2.876078287084027E-06"

$ws.Range("D10").Value = [double]"2.476675975991143E-39"

$ws.Range("D11").Value = [double]"8.090673646381473E-66"
$ws.Range("F11").Value = [double]"33.89786911010742"
